$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.277.92"
$ws.Range("E2").Value = "'  +0.23%  "
$ws.Range("D3").Value = "'3.501.14"
$ws.Range("E3").Value = "'  -0.50%  "
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("D5").Value = "'587.89"
$ws.Range("E5").Value = "'  +0.22%  "
$ws.Range("D6").Value = "'134.32"
$ws.Range("E6").Value = "'  +0.67%  "
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E8").Value = "'  -0.39%  "
$ws.Range("E9").Value = "'  +0.18%  "
$ws.Range("E10").Value = "'  +2.03%  "
$ws.Range("D11").Value = "'0.385"
$ws.Range("E11").Value = "'  +2.17%  "
$ws.Range("D12").Value = "'4.100.94"
$ws.Range("E12").Value = "'  -0.39%  "
$ws.Range("E13").Value = "'  +1.24%  "
$ws.Range("D14").Value = "'0.0000180"
$ws.Range("E14").Value = "'  +1.10%  "
$ws.Range("D15").Value = "'3.503.80"
$ws.Range("E15").Value = "'  -0.42%  "
$ws.Range("D16").Value = "'64.296.72"
$ws.Range("E16").Value = "'  +0.22%  "
$ws.Range("D17").Value = "'25.62"
$ws.Range("E17").Value = "'  -6.78%  "
$ws.Range("D18").Value = "'9.85"
$ws.Range("E18").Value = "'  +0.16%  "
$ws.Range("D19").Value = "'5.74"
$ws.Range("E19").Value = "'  +2.37%  "
$ws.Range("D20").Value = "'13.52"
$ws.Range("E20").Value = "'  -2.55%  "
$ws.Range("D21").Value = "'392.91"
$ws.Range("E21").Value = "'  +2.63%  "
$ws.Range("D22").Value = "'0.570"
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("D23").Value = "'3.641.61"
$ws.Range("E23").Value = "'  -0.50%  "
$ws.Range("D24").Value = "'74.57"
$ws.Range("E24").Value = "'  +0.85%  "
$ws.Range("D25").Value = "'1.01"
$ws.Range("E25").Value = "'  +1.10%  "
$ws.Range("B26").Value = "'PEPE"
$ws.Range("C26").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "'0.0000115"
$ws.Range("E26").Value = "'  +0.41%  "
$ws.Range("B27").Value = "'Binance-PegBSC-USD"
$ws.Range("C27").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "'  +0.03%  "
$ws.Range("B28").Value = "'RenderToken"
$ws.Range("C28").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.33"
$ws.Range("E28").Value = "'  -1.68%  "
$ws.Range("B29").Value = "'PancakeSwap"
$ws.Range("C29").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "'  +0.73%  "
$ws.Range("B30").Value = "'InternetComputer(DFINITY)"
$ws.Range("C30").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'8.23"
$ws.Range("E30").Value = "'  -2.72%  "
$ws.Range("B31").Value = "'Fetch.AI"
$ws.Range("C31").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.47"
$ws.Range("E31").Value = "'  -7.16%  "
$ws.Range("B32").Value = "'RenzoRestakedETH"
$ws.Range("C32").Value = "'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D32").Value = "'3.525.11"
$ws.Range("E32").Value = "'  -0.19%  "
$ws.Range("B33").Value = "'Kaspa"
$ws.Range("C33").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "'0.153"
$ws.Range("E33").Value = "'  +5.42%  "
$ws.Range("B34").Value = "'USDe"
$ws.Range("C34").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "'  +0.06%  "
$ws.Range("B35").Value = "'EthereumClassic"
$ws.Range("C35").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'23.41"
$ws.Range("E35").Value = "'  -0.61%  "
$ws.Range("B36").Value = "'NEARProtocol"
$ws.Range("C36").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.12"
$ws.Range("E36").Value = "'  -4.66%  "
$ws.Range("B37").Value = "'Aptos"
$ws.Range("C37").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'6.87"
$ws.Range("E37").Value = "'  -1.17%  "
$ws.Range("B38").Value = "'Monero"
$ws.Range("C38").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'167.62"
$ws.Range("E38").Value = "'  +4.53%  "
$ws.Range("E39").Value = "'  -1.27%  "
$ws.Range("B40").Value = "'Hedera"
$ws.Range("C40").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0778"
$ws.Range("E40").Value = "'  -0.77%  "
$ws.Range("B41").Value = "'Mantle"
$ws.Range("C41").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "'0.809"
$ws.Range("E41").Value = "'  -0.50%  "
$ws.Range("B42").Value = "'FirstDigitalUSD"
$ws.Range("C42").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "'  +0.07%  "
$ws.Range("B43").Value = "'EnergySwap"
$ws.Range("C43").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'25.31"
$ws.Range("E43").Value = "'  -4.82%  "
$ws.Range("B44").Value = "'Filecoin"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.39"
$ws.Range("E44").Value = "'  -0.33%  "
$ws.Range("B45").Value = "'Stacks"
$ws.Range("C45").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.65"
$ws.Range("E45").Value = "'  +2.72%  "
$ws.Range("B46").Value = "'ONDO"
$ws.Range("C46").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "'1.16"
$ws.Range("E46").Value = "'  -4.54%  "
$ws.Range("B47").Value = "'Cosmos"
$ws.Range("C47").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "'6.75"
$ws.Range("E47").Value = "'  -0.66%  "
$ws.Range("B48").Value = "'SuiNetwork"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D48").Value = "'0.888"
$ws.Range("E48").Value = "'  -2.30%  "
$ws.Range("B49").Value = "'Maker"
$ws.Range("C49").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'2.303.05"
$ws.Range("E49").Value = "'  -7.17%  "
$ws.Range("B50").Value = "'VeChain"
$ws.Range("C50").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0259"
$ws.Range("E50").Value = "'  -1.40%  "
$ws.Range("B51").Value = "'InjectiveProtocol"
$ws.Range("C51").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'21.14"
$ws.Range("E51").Value = "'  -1.50%  "
